$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 283, pushing the existing data (rows 283..339) down
# by one (to rows 284..340). This matches the diff: every existing row's
# data is shifted down one row, and a brand-new record appears at row 283.
$ws.Rows(283).Insert()

# The cells that were just pushed down to row 284 hold what used to be
# row 283's data. Re-use most of that row's field values (market id,
# market name, region, codreg, category id/name, variety, quality,
# volume, unit, origin, kg/units, classification) for the new row 283,
# since only the date and the three price columns + $/Kg actually change.
$ws.Cells.Item(283, 1).Value = $ws.Cells.Item(284, 1).Value()
$ws.Cells.Item(283, 2).Value = $ws.Cells.Item(284, 2).Value()
$ws.Cells.Item(283, 3).Value = $ws.Cells.Item(284, 3).Value()
$ws.Cells.Item(283, 4).Value = 45209
$ws.Cells.Item(283, 5).Value = $ws.Cells.Item(284, 5).Value()
$ws.Cells.Item(283, 6).Value = $ws.Cells.Item(284, 6).Value()
$ws.Cells.Item(283, 7).Value = $ws.Cells.Item(284, 7).Value()
$ws.Cells.Item(283, 8).Value = $ws.Cells.Item(284, 8).Value()
$ws.Cells.Item(283, 9).Value = $ws.Cells.Item(284, 9).Value()
$ws.Cells.Item(283, 10).Value = $ws.Cells.Item(284, 10).Value()
$ws.Cells.Item(283, 11).Value = 9000
$ws.Cells.Item(283, 12).Value = 9000
$ws.Cells.Item(283, 13).Value = 9000
$ws.Cells.Item(283, 14).Value = $ws.Cells.Item(284, 14).Value()
$ws.Cells.Item(283, 15).Value = $ws.Cells.Item(284, 15).Value()
$ws.Cells.Item(283, 16).Value = 750
$ws.Cells.Item(283, 17).Value = $ws.Cells.Item(284, 17).Value()
$ws.Cells.Item(283, 18).Value = $ws.Cells.Item(284, 18).Value()
